$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = [double]"1.139856140145738e-19"
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = [double]"4.828585634557355e-10"
$ws.Range("D2").Value = [double]"-3.219057405581146e-10"

$ws.Range("A3").Value = [double]"1.931434253822948e-10"
$ws.Range("B3").Value = [double]"4.828585634557358e-10"
$ws.Range("C3").Value = -1
$ws.Range("D3").Value = [double]"-1.22328035736726e-16"

$ws.Range("A4").Value = [double]"1.609528544852451e-10"
$ws.Range("B4").Value = [double]"3.219058626773073e-10"
$ws.Range("C4").Value = [double]"-6.717089639404731e-16"
$ws.Range("D4").Value = 1

$ws.Range("A5").Value = 1
$ws.Range("B5").Value = [double]"-2.590582384090349e-19"
$ws.Range("C5").Value = [double]"1.931434253822923e-10"
$ws.Range("D5").Value = [double]"-1.609528544852464e-10"
